# Scheduled-runner update: refresh cached Universalis market-price columns
# (currentAveragePrice / NQ / HQ, LevePriceNQ/HQ, LeveProfitNQ/HQ) on the
# per-job Leve-profit sheets. Values below are the new snapshot; blank
# ("") assignments reproduce cells the refresh dropped entirely (N/A
# profit rows), matching the upstream commit's cell-level diff exactly.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 0
$ws.Range("J32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("N32").Value = ""
$ws.Range("H40").Value = 5978.875
$ws.Range("I40").Value = 5542.7
$ws.Range("J40").Value = 6705.8335
$ws.Range("K40").Value = 5542.7
$ws.Range("L40").Value = 6705.8335
$ws.Range("M40").Value = -5367.7
$ws.Range("N40").Value = -7055.8335
$ws.Range("H125").Value = 2567.2727
$ws.Range("I125").Value = 2351.875
$ws.Range("K125").Value = 21166.875
$ws.Range("M125").Value = -18706.875
$ws.Range("H131").Value = 1814.5
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("N131").Value = ""
$ws.Range("H132").Value = 17528.732
$ws.Range("I132").Value = 17560.916
$ws.Range("K132").Value = 52682.74800000001
$ws.Range("M132").Value = -50152.74800000001

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 6904.5
$ws.Range("I2").Value = 4430.8184
$ws.Range("K2").Value = 4430.8184
$ws.Range("M2").Value = -4317.8184
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = ""
$ws.Range("N3").Value = ""
$ws.Range("H102").Value = 4740.5454
$ws.Range("I102").Value = 2768.25
$ws.Range("K102").Value = 2768.25
$ws.Range("M102").Value = -1146.25
$ws.Range("H110").Value = 2676.4736
$ws.Range("I110").Value = 2545.3845
$ws.Range("J110").Value = 2960.5
$ws.Range("K110").Value = 2545.3845
$ws.Range("L110").Value = 2960.5
$ws.Range("M110").Value = -500.3845000000001
$ws.Range("N110").Value = -7050.5
$ws.Range("H116").Value = 6904.5
$ws.Range("I116").Value = 4430.8184
$ws.Range("K116").Value = 4430.8184
$ws.Range("M116").Value = -2136.8184
$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").Value = ""

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 6904.5
$ws.Range("I3").Value = 4430.8184
$ws.Range("K3").Value = 4430.8184
$ws.Range("M3").Value = -4316.8184
$ws.Range("H86").Value = 3559.7778
$ws.Range("J86").Value = 5126.5
$ws.Range("L86").Value = 5126.5
$ws.Range("N86").Value = -7372.5
$ws.Range("H89").Value = 3559.7778
$ws.Range("J89").Value = 5126.5
$ws.Range("L89").Value = 25632.5
$ws.Range("N89").Value = -36864.5
$ws.Range("H94").Value = 864.25
$ws.Range("I94").Value = 819
$ws.Range("K94").Value = 819
$ws.Range("M94").Value = -368
$ws.Range("H99").Value = 2699.2
$ws.Range("I99").Value = 2699.2
$ws.Range("K99").Value = 2699.2
$ws.Range("M99").Value = -1201.2
$ws.Range("H105").Value = 4677
$ws.Range("I105").Value = 2010
$ws.Range("K105").Value = 2010
$ws.Range("M105").Value = -263
$ws.Range("H107").Value = 3671.7856
$ws.Range("I107").Value = 1244.4375
$ws.Range("J107").Value = 6908.25
$ws.Range("K107").Value = 1244.4375
$ws.Range("L107").Value = 6908.25
$ws.Range("M107").Value = 675.5625
$ws.Range("N107").Value = -10748.25

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 74.166664
$ws.Range("I7").Value = 75.25
$ws.Range("J7").Value = 72
$ws.Range("K7").Value = 75.25
$ws.Range("L7").Value = 72
$ws.Range("M7").Value = 37.75
$ws.Range("N7").Value = -298
$ws.Range("H31").Value = 7347
$ws.Range("J31").Value = 8746.799999999999
$ws.Range("L31").Value = 8746.799999999999
$ws.Range("N31").Value = -9336.799999999999
$ws.Range("H34").Value = 7347
$ws.Range("J34").Value = 8746.799999999999
$ws.Range("L34").Value = 8746.799999999999
$ws.Range("N34").Value = -9150.799999999999
$ws.Range("H47").Value = 24949.5
$ws.Range("I47").Value = 24949.5
$ws.Range("K47").Value = 24949.5
$ws.Range("M47").Value = -24383.5
$ws.Range("H105").Value = 2465.2
$ws.Range("I105").Value = 2465.2
$ws.Range("K105").Value = 2465.2
$ws.Range("M105").Value = -718.1999999999998
$ws.Range("H107").Value = 273.16666
$ws.Range("I107").Value = 327.66666
$ws.Range("J107").Value = 164.16667
$ws.Range("K107").Value = 327.66666
$ws.Range("L107").Value = 164.16667
$ws.Range("M107").Value = 1592.33334
$ws.Range("N107").Value = -4004.16667
$ws.Range("H134").Value = 2234
$ws.Range("I134").Value = 2234
$ws.Range("K134").Value = 6702
$ws.Range("M134").Value = -4167

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H70").Value = 500
$ws.Range("I70").Value = 500
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 1500
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = -1185
$ws.Range("N70").Value = ""
$ws.Range("H73").Value = 500
$ws.Range("I73").Value = 500
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 1500
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = -408
$ws.Range("N73").Value = ""

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 125055.875
$ws.Range("I2").Value = 250026.75
$ws.Range("J2").Value = 85
$ws.Range("K2").Value = 250026.75
$ws.Range("L2").Value = 85
$ws.Range("M2").Value = -249913.75
$ws.Range("N2").Value = -311
$ws.Range("H23").Value = 1145
$ws.Range("J23").Value = 1145
$ws.Range("L23").Value = 1145
$ws.Range("N23").Value = -1591
$ws.Range("H43").Value = 12349.75
$ws.Range("J43").Value = 11466.333
$ws.Range("L43").Value = 11466.333
$ws.Range("N43").Value = -11768.333
$ws.Range("H46").Value = 9657.666999999999
$ws.Range("I46").Value = 8974
$ws.Range("J46").Value = 9999.5
$ws.Range("K46").Value = 8974
$ws.Range("L46").Value = 9999.5
$ws.Range("M46").Value = -8818
$ws.Range("N46").Value = -10311.5
$ws.Range("H57").Value = 32726.75
$ws.Range("I57").Value = 19950
$ws.Range("J57").Value = 45503.5
$ws.Range("K57").Value = 19950
$ws.Range("L57").Value = 45503.5
$ws.Range("M57").Value = -19130
$ws.Range("N57").Value = -47143.5
$ws.Range("H97").Value = 927
$ws.Range("I97").Value = 737.2222
$ws.Range("K97").Value = 737.2222
$ws.Range("M97").Value = -241.2222
$ws.Range("H102").Value = 1258.9546
$ws.Range("I102").Value = 1258.9546
$ws.Range("K102").Value = 1258.9546
$ws.Range("M102").Value = 363.0454
$ws.Range("H134").Value = 125494
$ws.Range("J134").Value = 125494
$ws.Range("L134").Value = 376482
$ws.Range("N134").Value = -381552

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6102.231
$ws.Range("I7").Value = 4820.4443
$ws.Range("K7").Value = 4820.4443
$ws.Range("M7").Value = -4708.4443
$ws.Range("H32").Value = 1881
$ws.Range("I32").Value = 1881
$ws.Range("K32").Value = 1881
$ws.Range("M32").Value = -1564
$ws.Range("H40").Value = 4537.6665
$ws.Range("I40").Value = 2332.7
$ws.Range("J40").Value = 8947.6
$ws.Range("K40").Value = 2332.7
$ws.Range("L40").Value = 8947.6
$ws.Range("M40").Value = -2196.7
$ws.Range("N40").Value = -9219.6
$ws.Range("H61").Value = 3001.2942
$ws.Range("I61").Value = 1463.2307
$ws.Range("K61").Value = 1463.2307
$ws.Range("M61").Value = -1261.2307
$ws.Range("H93").Value = 1281.75
$ws.Range("I93").Value = 1281.75
$ws.Range("K93").Value = 1281.75
$ws.Range("M93").Value = -33.75
$ws.Range("H100").Value = 6058.9287
$ws.Range("I100").Value = 2385
$ws.Range("K100").Value = 2385
$ws.Range("M100").Value = -1844
$ws.Range("H113").Value = 3001.2942
$ws.Range("I113").Value = 1463.2307
$ws.Range("K113").Value = 1463.2307
$ws.Range("M113").Value = 706.7692999999999
$ws.Range("H122").Value = 4976.4
$ws.Range("I122").Value = 4976.875
$ws.Range("K122").Value = 14930.625
$ws.Range("M122").Value = -12480.625
$ws.Range("H126").Value = 6102.231
$ws.Range("I126").Value = 4820.4443
$ws.Range("K126").Value = 14461.3329
$ws.Range("M126").Value = -11991.3329

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 460.42856
$ws.Range("I107").Value = 353.83334
$ws.Range("J107").Value = 1100
$ws.Range("K107").Value = 1061.50002
$ws.Range("L107").Value = 3300
$ws.Range("M107").Value = 858.4999800000001
$ws.Range("N107").Value = -7140
$ws.Range("H126").Value = 3587.4285
$ws.Range("I126").Value = 1357
$ws.Range("J126").Value = 5817.857
$ws.Range("K126").Value = 4071
$ws.Range("L126").Value = 17453.571
$ws.Range("M126").Value = -1601
$ws.Range("N126").Value = -22393.571
$ws.Range("H132").Value = 3204.7273
$ws.Range("I132").Value = 1824.8334
$ws.Range("K132").Value = 5474.5002
$ws.Range("M132").Value = -2944.5002
$ws.Range("H138").Value = 95000
$ws.Range("J138").Value = 95000
$ws.Range("L138").Value = 95000
$ws.Range("N138").Value = -105280
